# Update "想去人数" (interested-count) figures in column F on the
# "展览" and "全部类型" sheets, per the latest data refresh.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value } for column F.
$updates = @{
    "展览"   = @{ 3 = 549; 6 = 502; 10 = 6705; 12 = 369; 13 = 2979; 15 = 336 }
    "全部类型" = @{ 5 = 549; 8 = 502; 13 = 6705; 16 = 369; 17 = 2979; 19 = 336 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
